$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Paragraph 1: merge the 9 runs of the first bio paragraph into a single,
# unformatted run; switch line spacing to 360/auto; add a first-line indent
# of 720 twips; and collapse the paragraph mark run properties down to just
# a language tag.
# ---------------------------------------------------------------------------
$para1Text = 'Sharon Jordan-Evans 在员工挽留和签约领域成绩骄人。她与 Beverly Kaye 合著了《华尔街日报》畅销书《Love ‘Em or Lose ‘Em:Getting Good People to Stay》，当前发行了第 4 版，被译成 20 多种语言。她的新作《Love It, Don’t Leave It:26 Ways to Get What You Want at Work》同样成为了《华尔街日报》畅销书，现已被翻译成 15 种语言。'

$para1Xml = '<w:p ' + $wNs + '><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t>' + $para1Text + '</w:t></w:r></w:p>'

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML($para1Xml)

# ---------------------------------------------------------------------------
# Paragraph 2: currently an empty paragraph. Give it the merged text of the
# second bio paragraph, switch on the first-line indent, drop the old
# paragraph-mark formatting, and plant the "_GoBack" bookmark at its start.
# ---------------------------------------------------------------------------
$para2Text = 'Sharon 经营着咨询公司 The Jordan Evans Group，主要培训高绩效高管并讲解敬业精神和员工挽留。她经常以企业辅导员和主题发言人的身份与美国证券交易所、波音、迪士尼、魔声公司、洛克希德公司和索尼等财富 500 强企业合作。她的个人网站：www.jeg.org'

$para2Xml = '<w:p ' + $wNs + '><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>' + $para2Text + '</w:t></w:r></w:p>'

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML($para2Xml)
